# IST price update 2025-12-20 17:54
#
# The tracker sheet keeps one column per price-check timestamp, newest
# first (right after the SKU name column). This run adds a brand new
# "latest check" column at B, pushing the previous columns one slot to
# the right (old B -> C, old C -> D) and stamping the new column's
# header with the current check time. Prices in the new column start
# out equal to the previous latest column's prices (most SKUs didn't
# move since the last check); a couple of rows had genuinely different
# B vs C prices before, and those differences simply ride along to the
# right as the columns shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B. This shifts all existing data,
# styles and formatting (old B -> C, old C -> D) and keeps the header
# cell style (s="1") intact on the shifted cells.
$ws.Columns("B").Insert()

# The inserted column doesn't inherit a width from its neighbours, so
# restore it to match columns C/D (internal OOXML width 21, which
# corresponds to this COM ColumnWidth value).
$ws.Columns("B").ColumnWidth = 20.17

# Figure out the full extent of the data (header + price rows).
$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

# New column B's price rows mirror column C (the just-shifted former
# column B) -- i.e. "no change since last check" for every SKU.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = $ws.Cells.Item($r, 3).Value()
}

# Stamp the header of the new latest-check column with this run's
# timestamp.
$ws.Range("B1").Value = "2025-12-20 23:18"
